$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 332.04166
$ws.Range("I33").Value = 358.5238
$ws.Range("K33").Value = 358.5238
$ws.Range("M33").Value = -129.5238
$ws.Range("H76").Value = 4546.8667
$ws.Range("I76").Value = 2922.5557
$ws.Range("J76").Value = 6983.3335
$ws.Range("K76").Value = 2922.5557
$ws.Range("L76").Value = 6983.3335
$ws.Range("M76").Value = -2607.5557
$ws.Range("N76").Value = -7613.3335
$ws.Range("H79").Value = 4546.8667
$ws.Range("I79").Value = 2922.5557
$ws.Range("J79").Value = 6983.3335
$ws.Range("K79").Value = 2922.5557
$ws.Range("L79").Value = 6983.3335
$ws.Range("M79").Value = -1830.5557
$ws.Range("N79").Value = -9167.333500000001
$ws.Range("H100").Value = 1626.5
$ws.Range("I100").Value = 1200
$ws.Range("J100").Value = 2337.3333
$ws.Range("K100").Value = 1200
$ws.Range("L100").Value = 2337.3333
$ws.Range("M100").Value = -659
$ws.Range("N100").Value = -3419.3333
$ws.Range("H107").Value = 7130.9414
$ws.Range("I107").Value = 7130.9414
$ws.Range("K107").Value = 7130.9414
$ws.Range("M107").Value = -5210.9414
$ws.Range("H116").Value = 3257.12
$ws.Range("I116").Value = 2430.7058
$ws.Range("J116").Value = 5013.25
$ws.Range("K116").Value = 2430.7058
$ws.Range("L116").Value = 5013.25
$ws.Range("M116").Value = 1011.2942
$ws.Range("N116").Value = -11897.25
$ws.Range("H137").Value = 6897395
$ws.Range("I137").Value = 762.0625
$ws.Range("J137").Value = 15385558
$ws.Range("K137").Value = 2286.1875
$ws.Range("L137").Value = 46156674
$ws.Range("M137").Value = 263.8125
$ws.Range("N137").Value = -46161774
$ws.Range("H138").Value = 1657.0444
$ws.Range("I138").Value = 1496.5897
$ws.Range("J138").Value = 2700
$ws.Range("K138").Value = 4489.7691
$ws.Range("L138").Value = 8100
$ws.Range("M138").Value = 650.2309000000005
$ws.Range("N138").Value = -18380
$ws.Range("H139").Value = 138096
$ws.Range("J139").Value = 138096
$ws.Range("L139").Value = 138096
$ws.Range("N139").Value = -148376
$ws.Range("H140").Value = 49700
$ws.Range("J140").Value = 49700
$ws.Range("L140").Value = 49700
$ws.Range("N140").Value = -60060
$ws.Range("H141").Value = 990.7727
$ws.Range("I141").Value = 990.7727
$ws.Range("K141").Value = 2972.3181
$ws.Range("M141").Value = 2207.6819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14287617
$ws.Range("I74").Value = 20001864
$ws.Range("J74").Value = 1999.4
$ws.Range("K74").Value = 20001864
$ws.Range("L74").Value = 1999.4
$ws.Range("M74").Value = -20000990
$ws.Range("N74").Value = -3747.4
$ws.Range("H77").Value = 14287617
$ws.Range("I77").Value = 20001864
$ws.Range("J77").Value = 1999.4
$ws.Range("K77").Value = 100009320
$ws.Range("L77").Value = 9997
$ws.Range("M77").Value = -100004952
$ws.Range("N77").Value = -18733
$ws.Range("H122").Value = 3291.5107
$ws.Range("I122").Value = 3457.2
$ws.Range("J122").Value = 2344.7144
$ws.Range("K122").Value = 10371.6
$ws.Range("L122").Value = 7034.1432
$ws.Range("M122").Value = -7921.599999999999
$ws.Range("N122").Value = -11934.1432
$ws.Range("H132").Value = 4547056.5
$ws.Range("I132").Value = 5556912.5
$ws.Range("J132").Value = 2704.2
$ws.Range("K132").Value = 16670737.5
$ws.Range("L132").Value = 8112.599999999999
$ws.Range("M132").Value = -16668207.5
$ws.Range("N132").Value = -13172.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 485.55554
$ws.Range("I22").Value = 485.55554
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 485.55554
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -312.55554
$ws.Range("N22").ClearContents()
$ws.Range("H134").Value = 2150.8462
$ws.Range("I134").Value = 1350.2549
$ws.Range("K134").Value = 4050.7647
$ws.Range("M134").Value = -1515.7647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 474.58334
$ws.Range("I22").Value = 299.375
$ws.Range("K22").Value = 299.375
$ws.Range("M22").Value = 50.625
$ws.Range("H31").Value = 4447641.5
$ws.Range("I31").Value = 3118.8225
$ws.Range("J31").Value = 25644594
$ws.Range("K31").Value = 3118.8225
$ws.Range("L31").Value = 25644594
$ws.Range("M31").Value = -2823.8225
$ws.Range("N31").Value = -25645184
$ws.Range("H34").Value = 4447641.5
$ws.Range("I34").Value = 3118.8225
$ws.Range("J34").Value = 25644594
$ws.Range("K34").Value = 3118.8225
$ws.Range("L34").Value = 25644594
$ws.Range("M34").Value = -2916.8225
$ws.Range("N34").Value = -25644998
$ws.Range("H75").Value = 50245
$ws.Range("J75").Value = 50245
$ws.Range("L75").Value = 50245
$ws.Range("N75").Value = -52241
$ws.Range("H78").Value = 50245
$ws.Range("J78").Value = 50245
$ws.Range("L78").Value = 150735
$ws.Range("N78").Value = -160719
$ws.Range("H125").Value = 22000
$ws.Range("J125").Value = 22000
$ws.Range("L125").Value = 22000
$ws.Range("N125").Value = -26920
$ws.Range("H132").Value = 11629415
$ws.Range("I132").Value = 12196364
$ws.Range("K132").Value = 36589092
$ws.Range("M132").Value = -36586562
$ws.Range("H134").Value = 1397.3265
$ws.Range("I134").Value = 1268.8914
$ws.Range("J134").Value = 3366.6667
$ws.Range("K134").Value = 3806.6742
$ws.Range("L134").Value = 10100.0001
$ws.Range("M134").Value = -1271.6742
$ws.Range("N134").Value = -15170.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5842
$ws.Range("I3").Value = 3850.4167
$ws.Range("J3").Value = 6881.087
$ws.Range("K3").Value = 11551.2501
$ws.Range("L3").Value = 20643.261
$ws.Range("M3").Value = -11439.2501
$ws.Range("N3").Value = -20867.261
$ws.Range("H129").Value = 4507.7334
$ws.Range("I129").Value = 3375
$ws.Range("J129").Value = 5262.8887
$ws.Range("K129").Value = 10125
$ws.Range("L129").Value = 15788.6661
$ws.Range("M129").Value = -5125
$ws.Range("N129").Value = -25788.6661
$ws.Range("H133").Value = 7140
$ws.Range("I133").Value = 3600
$ws.Range("K133").Value = 10800
$ws.Range("M133").Value = -5740
$ws.Range("H137").Value = 5994
$ws.Range("I137").Value = 3021.75
$ws.Range("J137").Value = 8092.0586
$ws.Range("K137").Value = 9065.25
$ws.Range("L137").Value = 24276.1758
$ws.Range("M137").Value = -3965.25
$ws.Range("N137").Value = -34476.1758
$ws.Range("H139").Value = 2671.9167
$ws.Range("I139").Value = 1220
$ws.Range("K139").Value = 3660
$ws.Range("M139").Value = 1480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 40780
$ws.Range("J130").Value = 40780
$ws.Range("L130").Value = 40780
$ws.Range("N130").Value = -50820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1601.95
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705
$ws.Range("H27").Value = 1601.95
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893
$ws.Range("H68").Value = 1337.1428
$ws.Range("I68").Value = 1540
$ws.Range("J68").Value = 830
$ws.Range("K68").Value = 1540
$ws.Range("L68").Value = 830
$ws.Range("M68").Value = -791
$ws.Range("N68").Value = -2328
$ws.Range("H71").Value = 1337.1428
$ws.Range("I71").Value = 1540
$ws.Range("J71").Value = 830
$ws.Range("K71").Value = 7700
$ws.Range("L71").Value = 4150
$ws.Range("M71").Value = -3956
$ws.Range("N71").Value = -11638
$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040
$ws.Range("H139").Value = 55716.668
$ws.Range("J139").Value = 55716.668
$ws.Range("L139").Value = 55716.668
$ws.Range("N139").Value = -65996.66800000001

Write-Host "Applied all changes"
